$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 897.36365
$ws.Range("I40").Value = 874.4286
$ws.Range("J40").Value = 937.5
$ws.Range("K40").Value = 874.4286
$ws.Range("L40").Value = 937.5
$ws.Range("M40").Value = -699.4286
$ws.Range("N40").Value = -1287.5
$ws.Range("H104").Value = 367.83334
$ws.Range("J104").Value = 942
$ws.Range("L104").Value = 2826
$ws.Range("N104").Value = -6320
$ws.Range("H129").Value = 552610.1
$ws.Range("I129").Value = 100539
$ws.Range("J129").Value = 1004681.2
$ws.Range("K129").Value = 301617
$ws.Range("L129").Value = 3014043.6
$ws.Range("M129").Value = -296617
$ws.Range("N129").Value = -3024043.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10877000
$ws.Range("I32").Value = 7464.095
$ws.Range("J32").Value = 125007130
$ws.Range("K32").Value = 7464.095
$ws.Range("L32").Value = 125007130
$ws.Range("M32").Value = -7177.095
$ws.Range("N32").Value = -125007704
$ws.Range("H45").Value = 1556.2222
$ws.Range("I45").Value = 1445.7778
$ws.Range("J45").Value = 1666.6666
$ws.Range("K45").Value = 1445.7778
$ws.Range("L45").Value = 1666.6666
$ws.Range("M45").Value = -1068.7778
$ws.Range("N45").Value = -2420.6666
$ws.Range("H110").Value = 1973.1111
$ws.Range("I110").Value = 2094.75
$ws.Range("J110").Value = 1000
$ws.Range("K110").Value = 2094.75
$ws.Range("L110").Value = 1000
$ws.Range("M110").Value = -49.75
$ws.Range("N110").Value = -5090

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1899.0825
$ws.Range("I86").Value = 1910.5685
$ws.Range("J86").Value = 1353.5
$ws.Range("K86").Value = 1910.5685
$ws.Range("L86").Value = 1353.5
$ws.Range("M86").Value = -787.5685000000001
$ws.Range("N86").Value = -3599.5
$ws.Range("H89").Value = 1899.0825
$ws.Range("I89").Value = 1910.5685
$ws.Range("J89").Value = 1353.5
$ws.Range("K89").Value = 9552.842500000001
$ws.Range("L89").Value = 6767.5
$ws.Range("M89").Value = -3936.842500000001
$ws.Range("N89").Value = -17999.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 50.923077
$ws.Range("I7").Value = 42.2
$ws.Range("J7").Value = 80
$ws.Range("K7").Value = 42.2
$ws.Range("L7").Value = 80
$ws.Range("M7").Value = 70.8
$ws.Range("N7").Value = -306
$ws.Range("H31").Value = 870981.75
$ws.Range("I31").Value = 1220.2122
$ws.Range("J31").Value = 1606933.8
$ws.Range("K31").Value = 1220.2122
$ws.Range("L31").Value = 1606933.8
$ws.Range("M31").Value = -925.2121999999999
$ws.Range("N31").Value = -1607523.8
$ws.Range("H34").Value = 870981.75
$ws.Range("I34").Value = 1220.2122
$ws.Range("J34").Value = 1606933.8
$ws.Range("K34").Value = 1220.2122
$ws.Range("L34").Value = 1606933.8
$ws.Range("M34").Value = -1018.2122
$ws.Range("N34").Value = -1607337.8
$ws.Range("H59").Value = 23057.143
$ws.Range("J59").Value = 23057.143
$ws.Range("L59").Value = 23057.143
$ws.Range("N59").Value = -25347.143
$ws.Range("H62").Value = 2560
$ws.Range("I62").Value = 2463.125
$ws.Range("K62").Value = 2463.125
$ws.Range("M62").Value = -1839.125
$ws.Range("H65").Value = 2560
$ws.Range("I65").Value = 2463.125
$ws.Range("K65").Value = 12315.625
$ws.Range("M65").Value = -9195.625
$ws.Range("H68").Value = 24319
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 24319
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 24319
$ws.Range("M68").Value = $null
$ws.Range("N68").Value = -25817
$ws.Range("H71").Value = 24319
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 24319
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 72957
$ws.Range("M71").Value = $null
$ws.Range("N71").Value = -80445
$ws.Range("H74").Value = 26634.75
$ws.Range("I74").Value = 20001
$ws.Range("J74").Value = 28846
$ws.Range("K74").Value = 20001
$ws.Range("L74").Value = 28846
$ws.Range("M74").Value = -19127
$ws.Range("N74").Value = -30594
$ws.Range("H77").Value = 26634.75
$ws.Range("I77").Value = 20001
$ws.Range("J77").Value = 28846
$ws.Range("K77").Value = 60003
$ws.Range("L77").Value = 86538
$ws.Range("M77").Value = -55635
$ws.Range("N77").Value = -95274
$ws.Range("H102").Value = 21142.857
$ws.Range("J102").Value = 21142.857
$ws.Range("L102").Value = 21142.857
$ws.Range("N102").Value = -26010.857
$ws.Range("H104").Value = 28000
$ws.Range("J104").Value = 28000
$ws.Range("L104").Value = 28000
$ws.Range("N104").Value = -33242
$ws.Range("H109").Value = 40285
$ws.Range("J109").Value = 40285
$ws.Range("L109").Value = 40285
$ws.Range("N109").Value = -42365

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 4263.5
$ws.Range("I94").Value = 500
$ws.Range("J94").Value = 5016.2
$ws.Range("K94").Value = 1500
$ws.Range("L94").Value = 15048.6
$ws.Range("M94").Value = -824
$ws.Range("N94").Value = -16400.6
$ws.Range("H131").Value = 8638909
$ws.Range("I131").Value = 50000290
$ws.Range("J131").Value = 21954.812
$ws.Range("K131").Value = 150000870
$ws.Range("L131").Value = 65864.436
$ws.Range("M131").Value = -149995830
$ws.Range("N131").Value = -75944.436

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 25000240
$ws.Range("I55").Value = 31250176
$ws.Range("J55").Value = 501
$ws.Range("K55").Value = 31250176
$ws.Range("L55").Value = 501
$ws.Range("M55").Value = -31250003
$ws.Range("N55").Value = -847
$ws.Range("H68").Value = 2418.3333
$ws.Range("I68").Value = 1500
$ws.Range("J68").Value = 2501.818
$ws.Range("K68").Value = 1500
$ws.Range("L68").Value = 2501.818
$ws.Range("M68").Value = -751
$ws.Range("N68").Value = -3999.818
$ws.Range("H71").Value = 2418.3333
$ws.Range("I71").Value = 1500
$ws.Range("J71").Value = 2501.818
$ws.Range("K71").Value = 7500
$ws.Range("L71").Value = 12509.09
$ws.Range("M71").Value = -3756
$ws.Range("N71").Value = -19997.09
